$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "64.854.72"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "3.396.50"
$ws.Range("E3").Value = "  -3.10%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.19%  "
Set-TextValue "D5" "578.87"
$ws.Range("E5").Value = "  -3.18%  "
Set-TextValue "D6" "134.60"
$ws.Range("E6").Value = "  -5.49%  "
Set-TextValue "D7" "1.00"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.385.73"
$ws.Range("E8").Value = "  -3.51%  "
Set-TextValue "D9" "0.491"
$ws.Range("E9").Value = "  -2.28%  "
Set-TextValue "D10" "0.120"
$ws.Range("E10").Value = "  -9.61%  "
Set-TextValue "D11" "7.04"
$ws.Range("E11").Value = "  -9.76%  "
Set-TextValue "D12" "0.370"
$ws.Range("E12").Value = "  -7.55%  "
$ws.Range("D13").Value = "3.977.74"
$ws.Range("E13").Value = "  -3.20%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D14" "0.0000176"
$ws.Range("E14").Value = "  -10.71%  "
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D15" "0.115"
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.399.13"
$ws.Range("E16").Value = "  -2.56%  "
$ws.Range("D17").Value = "64.823.54"
$ws.Range("E17").Value = "  -0.95%  "
Set-TextValue "D18" "25.88"
$ws.Range("E18").Value = "  -8.48%  "
Set-TextValue "D19" "9.47"
$ws.Range("E19").Value = "  -13.88%  "
Set-TextValue "D20" "5.76"
$ws.Range("E20").Value = "  -6.34%  "
$ws.Range("E21").Value = "  -5.88%  "
Set-TextValue "D22" "378.89"
$ws.Range("E22").Value = "  -8.69%  "
Set-TextValue "D23" "0.546"
$ws.Range("E23").Value = "  -7.90%  "
$ws.Range("E24").Value = "  +0.09%  "
Set-TextValue "D25" "71.71"
$ws.Range("E25").Value = "  -7.21%  "
$ws.Range("D26").Value = "3.533.87"
$ws.Range("E26").Value = "  -3.29%  "
Set-TextValue "D27" "0.0000103"
$ws.Range("E27").Value = "  -9.47%  "
Set-TextValue "D28" "0.997"
$ws.Range("E28").Value = "  +0.85%  "
Set-TextValue "D29" "6.94"
$ws.Range("E29").Value = "  -9.21%  "
Set-TextValue "D30" "2.17"
$ws.Range("E30").Value = "  -10.14%  "
Set-TextValue "D31" "7.90"
$ws.Range("E31").Value = "  -9.79%  "
$ws.Range("D32").Value = "3.411.70"
$ws.Range("E32").Value = "  -2.90%  "
$ws.Range("E34").Value = "  -7.10%  "
Set-TextValue "D35" "22.73"
$ws.Range("E35").Value = "  -5.92%  "
Set-TextValue "D36" "169.52"
$ws.Range("E36").Value = "  -1.64%  "
Set-TextValue "D37" "6.60"
$ws.Range("E37").Value = "  -11.50%  "
Set-TextValue "D38" "1.13"
$ws.Range("E38").Value = "  -11.82%  "
$ws.Range("E39").Value = "  -7.67%  "
Set-TextValue "D40" "4.60"
$ws.Range("E40").Value = "  -11.95%  "
Set-TextValue "D41" "0.0744"
$ws.Range("E41").Value = "  -8.33%  "
Set-TextValue "D42" "0.805"
$ws.Range("E42").Value = "  -5.39%  "
Set-TextValue "D43" "42.73"
$ws.Range("E43").Value = "  -5.26%  "
Set-TextValue "D44" "0.996"
$ws.Range("E44").Value = "  -0.53%  "
Set-TextValue "D45" "4.29"
$ws.Range("E45").Value = "  -14.58%  "
$ws.Range("E46").Value = "  -10.42%  "
Set-TextValue "D47" "1.08"
$ws.Range("E47").Value = "  +1.14%  "
Set-TextValue "D48" "21.84"
$ws.Range("E48").Value = "  -4.41%  "
Set-TextValue "D49" "6.40"
$ws.Range("E49").Value = "  -8.30%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D50" "2.00"
$ws.Range("E50").Value = "  -14.21%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.145.86"
$ws.Range("E51").Value = "  -8.48%  "
